$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the empty paragraph right after the title so the "Purpose:"
#    paragraph immediately follows the title paragraph.
# ---------------------------------------------------------------------------
$d.Paragraphs(2).Range.Delete()

# ---------------------------------------------------------------------------
# 2) Split "Installation:" into three runs with identical (underlined)
#    formatting: "Install" + "ation/Removal Instructions" + ":" -- i.e.
#    the middle substring "ation" is replaced by the longer text, while the
#    untouched "Install" prefix and ":" suffix remain / become their own
#    runs. Toggling a character property on then back off on a range is
#    enough to force a run break at that boundary even though the final
#    formatting ends up identical to its neighbour.
# ---------------------------------------------------------------------------
$installFind = $d.Content
$null = $installFind.Find.Execute("Installation:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$installStart = $installFind.Start

$prefixLen = "Install".Length
$midOld = "ation"
$midNew = "ation/Removal Instructions"

$mid = $d.Range($installStart + $prefixLen, $installStart + $prefixLen + $midOld.Length)
$mid.Text = $midNew

$installRun = $d.Range($installStart, $installStart + $prefixLen)
$installRun.Bold = 1
$installRun.Bold = 0

$midRun = $d.Range($installStart + $prefixLen, $installStart + $prefixLen + $midNew.Length)
$midRun.Bold = 1
$midRun.Bold = 0

# ---------------------------------------------------------------------------
# 3) Append the uninstall instructions to the end of the "Double-click..."
#    paragraph (just before its paragraph mark) as three additional runs.
# ---------------------------------------------------------------------------
$dblClickFind = $d.Content
$null = $dblClickFind.Find.Execute("Double-click the setup.exe file", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$installPara = $dblClickFind.Paragraphs(1)
$r = $installPara.Range
$insertStart = $r.End - 1

$t1 = 'Uninstall like a normal program: Add/Remove Program > select "IMS" > Uninstall'
$t2 = ' > OK'
$t3 = '.'

$r.InsertAfter($t1 + $t2 + $t3)

$run1 = $d.Range($insertStart, $insertStart + $t1.Length)
$run1.Bold = 1
$run1.Bold = 0

$run2 = $d.Range($insertStart + $t1.Length, $insertStart + $t1.Length + $t2.Length)
$run2.Bold = 1
$run2.Bold = 0

$run3 = $d.Range($insertStart + $t1.Length + $t2.Length, $insertStart + $t1.Length + $t2.Length + $t3.Length)
$run3.Bold = 1
$run3.Bold = 0
